# The edit deletes the recorded match result for "Teresa-Leticia" vs
# "Las barbas" (group "Mediocre medio") on the "resultados" sheet, row 5
# (A5:F5 -> blank, keeping E5's existing cell style). That, in turn, zeroes
# out the stats that had already been tallied for that match in the two
# classification sheets ("clasificacion_auto" and "clasificacion").

$wb = $excel.ActiveWorkbook

# --- resultados: clear the 5th match row, keep the selection on it -------
$wsResultados = $wb.Worksheets.Item("resultados")
[void]$wsResultados.Range("A5:F5").ClearContents()
[void]$wsResultados.Range("A5:XFD5").Select()

# --- clasificacion_auto: zero out the Teresa-Leticia / Las barbas stats --
$wsAuto = $wb.Worksheets.Item("clasificacion_auto")

# Teresa-Leticia (row 14): PUNTOS, PJ, PG, SG -> 0 (PE/PP/SP already 0)
$wsAuto.Range("D14").Value = 0
$wsAuto.Range("E14").Value = 0
$wsAuto.Range("F14").Value = 0
$wsAuto.Range("I14").Value = 0

# Las barbas (row 15): PJ, PP, SP -> 0 (rest already 0)
$wsAuto.Range("E15").Value = 0
$wsAuto.Range("H15").Value = 0
$wsAuto.Range("J15").Value = 0

# --- clasificacion: same reset, different row positions ------------------
$wsClasif = $wb.Worksheets.Item("clasificacion")

# Teresa-Leticia (row 7)
$wsClasif.Range("D7").Value = 0
$wsClasif.Range("E7").Value = 0
$wsClasif.Range("F7").Value = 0
$wsClasif.Range("I7").Value = 0

# Las barbas (row 8)
$wsClasif.Range("E8").Value = 0
$wsClasif.Range("H8").Value = 0
$wsClasif.Range("J8").Value = 0
